$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1975.7693
$ws.Range("I15").Value = 1975.7693
$ws.Range("K15").Value = 5927.3079
$ws.Range("M15").Value = -5758.3079

$ws.Range("H17").Value = 5433.3335
$ws.Range("J17").Value = 6000
$ws.Range("L17").Value = 18000
$ws.Range("N17").Value = -18336

$ws.Range("H18").Value = 258.27274
$ws.Range("I18").Value = 258.27274
$ws.Range("K18").Value = 258.27274
$ws.Range("M18").Value = 25.72726

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H28").Value = 927.6070999999999
$ws.Range("I28").Value = 605.7917
$ws.Range("K28").Value = 605.7917
$ws.Range("M28").Value = -120.7917

$ws.Range("H32").Value = 3649.6667
$ws.Range("I32").Value = 3724.5
$ws.Range("J32").Value = 3500
$ws.Range("K32").Value = 3724.5
$ws.Range("L32").Value = 3500
$ws.Range("M32").Value = -3398.5
$ws.Range("N32").Value = -4152

$ws.Range("H33").Value = 3380888.5
$ws.Range("I33").Value = 6756929.5
$ws.Range("K33").Value = 6756929.5
$ws.Range("M33").Value = -6756700.5

$ws.Range("H43").Value = 2720
$ws.Range("I43").Value = 2334.1904
$ws.Range("J43").Value = 4745.5
$ws.Range("K43").Value = 2334.1904
$ws.Range("L43").Value = 4745.5
$ws.Range("M43").Value = -2265.1904
$ws.Range("N43").Value = -4883.5

$ws.Range("H51").Value = 6562.4585
$ws.Range("I51").Value = 4500
$ws.Range("K51").Value = 4500
$ws.Range("M51").Value = -4016

$ws.Range("H64").Value = 6399.4287
$ws.Range("J64").Value = 7907.5386
$ws.Range("L64").Value = 7907.5386
$ws.Range("N64").Value = -8403.5386

$ws.Range("H67").Value = 6399.4287
$ws.Range("J67").Value = 7907.5386
$ws.Range("L67").Value = 7907.5386
$ws.Range("N67").Value = -9623.5386

$ws.Range("H74").Value = 7934
$ws.Range("I74").Value = 7934
$ws.Range("K74").Value = 7934
$ws.Range("M74").Value = -6998

$ws.Range("H76").Value = 3683.9285
$ws.Range("J76").Value = 3669.7144
$ws.Range("L76").Value = 3669.7144
$ws.Range("N76").Value = -4299.7144

$ws.Range("H77").Value = 7934
$ws.Range("I77").Value = 7934
$ws.Range("K77").Value = 39670
$ws.Range("M77").Value = -34990

$ws.Range("H79").Value = 3683.9285
$ws.Range("J79").Value = 3669.7144
$ws.Range("L79").Value = 3669.7144
$ws.Range("N79").Value = -5853.7144

$ws.Range("H80").Value = 359349.28
$ws.Range("I80").Value = 1885.7333
$ws.Range("J80").Value = 771807.25
$ws.Range("K80").Value = 5657.199900000001
$ws.Range("L80").Value = 2315421.75
$ws.Range("M80").Value = -4659.199900000001
$ws.Range("N80").Value = -2317417.75

$ws.Range("H83").Value = 359349.28
$ws.Range("I83").Value = 1885.7333
$ws.Range("J83").Value = 771807.25
$ws.Range("K83").Value = 16971.5997
$ws.Range("L83").Value = 6946265.25
$ws.Range("M83").Value = -11979.5997
$ws.Range("N83").Value = -6956249.25

$ws.Range("H87").Value = 81245
$ws.Range("J87").Value = 83328
$ws.Range("L87").Value = 83328
$ws.Range("N87").Value = -85824

$ws.Range("H90").Value = 81245
$ws.Range("J90").Value = 83328
$ws.Range("L90").Value = 249984
$ws.Range("N90").Value = -262464

$ws.Range("H94").Value = 1748.8462
$ws.Range("I94").Value = 1436.25
$ws.Range("K94").Value = 1436.25
$ws.Range("M94").Value = -985.25

$ws.Range("H96").Value = 486.29413
$ws.Range("I96").Value = 344.46667
$ws.Range("K96").Value = 1033.40001
$ws.Range("M96").Value = 339.5999899999999

$ws.Range("H98").Value = 2285.7778
$ws.Range("I98").Value = 2214.353
$ws.Range("J98").Value = 3500
$ws.Range("K98").Value = 2214.353
$ws.Range("L98").Value = 3500
$ws.Range("M98").Value = -716.3530000000001
$ws.Range("N98").Value = -6496

$ws.Range("H101").Value = 1093.75
$ws.Range("I101").Value = 1235.6364
$ws.Range("K101").Value = 3706.9092
$ws.Range("M101").Value = -2084.9092

$ws.Range("H107").Value = 12657.608
$ws.Range("I107").Value = 7160.625
$ws.Range("J107").Value = 25222.143
$ws.Range("K107").Value = 7160.625
$ws.Range("L107").Value = 25222.143
$ws.Range("M107").Value = -5240.625
$ws.Range("N107").Value = -29062.143

$ws.Range("H112").Value = 1789.3846
$ws.Range("I112").Value = 1397.6666
$ws.Range("J112").Value = 1906.9
$ws.Range("K112").Value = 4192.9998
$ws.Range("L112").Value = 5720.700000000001
$ws.Range("M112").Value = -3084.9998
$ws.Range("N112").Value = -7936.700000000001

$ws.Range("H116").Value = 8813.117
$ws.Range("I116").Value = 7653.0835
$ws.Range("J116").Value = 11597.2
$ws.Range("K116").Value = 7653.0835
$ws.Range("L116").Value = 11597.2
$ws.Range("M116").Value = -4211.0835
$ws.Range("N116").Value = -18481.2

$ws.Range("H121").Value = 3766.4443
$ws.Range("J121").Value = 4049.875
$ws.Range("L121").Value = 12149.625
$ws.Range("N121").Value = -15643.625

$ws.Range("H122").Value = 2285.7778
$ws.Range("I122").Value = 2214.353
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 6643.059
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -4193.059
$ws.Range("N122").Value = -15400

$ws.Range("H123").Value = 49999.332
$ws.Range("J123").Value = 49999.332
$ws.Range("L123").Value = 49999.332
$ws.Range("N123").Value = -59799.332

$ws.Range("H132").Value = 2237.0732
$ws.Range("I132").Value = 2239
$ws.Range("K132").Value = 6717
$ws.Range("M132").Value = -4187

$ws.Range("H135").Value = 4438.4375
$ws.Range("I135").Value = 4860.6924
$ws.Range("K135").Value = 43746.2316
$ws.Range("M135").Value = -41211.2316

$ws.Range("H137").Value = 12327
$ws.Range("I137").Value = 4497.9414
$ws.Range("K137").Value = 13493.8242
$ws.Range("M137").Value = -10943.8242

$ws.Range("H138").Value = 1971.075
$ws.Range("I138").Value = 1721.4572
$ws.Range("J138").Value = 3718.4
$ws.Range("K138").Value = 5164.3716
$ws.Range("L138").Value = 11155.2
$ws.Range("M138").Value = -24.3716000000004
$ws.Range("N138").Value = -21435.2

$ws.Range("H141").Value = 2117.6296
$ws.Range("I141").Value = 2022.2273
$ws.Range("J141").Value = 2537.4
$ws.Range("K141").Value = 6066.6819
$ws.Range("L141").Value = 7612.200000000001
$ws.Range("M141").Value = -886.6818999999996
$ws.Range("N141").Value = -17972.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6368.8047
$ws.Range("I32").Value = 5669.769
$ws.Range("K32").Value = 5669.769
$ws.Range("M32").Value = -5382.769

$ws.Range("H45").Value = 6071.16
$ws.Range("J45").Value = 1706.25
$ws.Range("L45").Value = 1706.25
$ws.Range("N45").Value = -2460.25

$ws.Range("H61").Value = 13745.435
$ws.Range("I61").Value = 9557.666999999999
$ws.Range("K61").Value = 9557.666999999999
$ws.Range("M61").Value = -9345.666999999999

$ws.Range("H63").Value = 3322.9285
$ws.Range("I63").Value = 2020.625
$ws.Range("K63").Value = 2020.625
$ws.Range("M63").Value = -1334.625

$ws.Range("H66").Value = 3322.9285
$ws.Range("I66").Value = 2020.625
$ws.Range("K66").Value = 10103.125
$ws.Range("M66").Value = -6671.125

$ws.Range("H74").Value = 21839.264
$ws.Range("I74").Value = 25769.363
$ws.Range("J74").Value = 16435.375
$ws.Range("K74").Value = 25769.363
$ws.Range("L74").Value = 16435.375
$ws.Range("M74").Value = -24895.363
$ws.Range("N74").Value = -18183.375

$ws.Range("H77").Value = 21839.264
$ws.Range("I77").Value = 25769.363
$ws.Range("J77").Value = 16435.375
$ws.Range("K77").Value = 128846.815
$ws.Range("L77").Value = 82176.875
$ws.Range("M77").Value = -124478.815
$ws.Range("N77").Value = -90912.875

$ws.Range("H86").Value = 20000
$ws.Range("I86").Value = 20000
$ws.Range("K86").Value = 20000
$ws.Range("M86").Value = -18814

$ws.Range("H89").Value = 20000
$ws.Range("I89").Value = 20000
$ws.Range("K89").Value = 60000
$ws.Range("M89").Value = -54072

$ws.Range("H93").Value = 40448
$ws.Range("J93").Value = 40448
$ws.Range("L93").Value = 40448
$ws.Range("N93").Value = -45440

$ws.Range("H97").Value = 1324.6897
$ws.Range("I97").Value = 1204.625
$ws.Range("K97").Value = 1204.625
$ws.Range("M97").Value = -708.625

$ws.Range("H102").Value = 1943.7333
$ws.Range("I102").Value = 1975.4286
$ws.Range("K102").Value = 1975.4286
$ws.Range("M102").Value = -353.4286

$ws.Range("H119").Value = 57000
$ws.Range("J119").Value = 57000
$ws.Range("L119").Value = 57000
$ws.Range("N119").Value = -66676

$ws.Range("H132").Value = 1582.826
$ws.Range("I132").Value = 1443.9445
$ws.Range("K132").Value = 4331.833500000001
$ws.Range("M132").Value = -1801.833500000001

$ws.Range("H136").Value = 13745.435
$ws.Range("I136").Value = 9557.666999999999
$ws.Range("K136").Value = 28673.001
$ws.Range("M136").Value = -26123.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11451.911
$ws.Range("I20").Value = 11283.027
$ws.Range("K20").Value = 11283.027
$ws.Range("M20").Value = -11036.027

$ws.Range("H22").Value = 332.05
$ws.Range("I22").Value = 327.8889
$ws.Range("K22").Value = 327.8889
$ws.Range("M22").Value = -154.8889

$ws.Range("H80").Value = 307.0625
$ws.Range("I80").Value = 223.27272
$ws.Range("J80").Value = 491.4
$ws.Range("K80").Value = 223.27272
$ws.Range("L80").Value = 491.4
$ws.Range("M80").Value = 774.7272800000001
$ws.Range("N80").Value = -2487.4

$ws.Range("H83").Value = 307.0625
$ws.Range("I83").Value = 223.27272
$ws.Range("J83").Value = 491.4
$ws.Range("K83").Value = 1116.3636
$ws.Range("L83").Value = 2457
$ws.Range("M83").Value = 3875.6364
$ws.Range("N83").Value = -12441

$ws.Range("H86").Value = 336314.9
$ws.Range("I86").Value = 590089.8
$ws.Range("K86").Value = 590089.8
$ws.Range("M86").Value = -588966.8

$ws.Range("H89").Value = 336314.9
$ws.Range("I89").Value = 590089.8
$ws.Range("K89").Value = 2950449
$ws.Range("M89").Value = -2944833

$ws.Range("H94").Value = 1482.2174
$ws.Range("I94").Value = 1192.05
$ws.Range("J94").Value = 3416.6667
$ws.Range("K94").Value = 1192.05
$ws.Range("L94").Value = 3416.6667
$ws.Range("M94").Value = -741.05
$ws.Range("N94").Value = -4318.6667

$ws.Range("H99").Value = 6452.7812
$ws.Range("I99").Value = 6499.6553
$ws.Range("J99").Value = 5999.6665
$ws.Range("K99").Value = 6499.6553
$ws.Range("L99").Value = 5999.6665
$ws.Range("M99").Value = -5001.6553
$ws.Range("N99").Value = -8995.666499999999

$ws.Range("H134").Value = 16430.521
$ws.Range("I134").Value = 9438
$ws.Range("K134").Value = 28314
$ws.Range("M134").Value = -25779

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 67999.5
$ws.Range("J9").Value = 67999.5
$ws.Range("L9").Value = 67999.5
$ws.Range("N9").Value = -68335.5

$ws.Range("H22").Value = 472.82352
$ws.Range("I22").Value = 174.23077
$ws.Range("K22").Value = 174.23077
$ws.Range("M22").Value = 175.76923

$ws.Range("H31").Value = 2575.1428
$ws.Range("I31").Value = 1423.091
$ws.Range("J31").Value = 3842.4
$ws.Range("K31").Value = 1423.091
$ws.Range("L31").Value = 3842.4
$ws.Range("M31").Value = -1128.091
$ws.Range("N31").Value = -4432.4

$ws.Range("H34").Value = 2575.1428
$ws.Range("I34").Value = 1423.091
$ws.Range("J34").Value = 3842.4
$ws.Range("K34").Value = 1423.091
$ws.Range("L34").Value = 3842.4
$ws.Range("M34").Value = -1221.091
$ws.Range("N34").Value = -4246.4

$ws.Range("H43").Value = 13861.2
$ws.Range("J43").Value = 13861.2
$ws.Range("L43").Value = 13861.2
$ws.Range("N43").Value = -14229.2

$ws.Range("H58").Value = 3838.75
$ws.Range("I58").Value = 2201.28
$ws.Range("J58").Value = 5993.316
$ws.Range("K58").Value = 2201.28
$ws.Range("L58").Value = 5993.316
$ws.Range("M58").Value = -1998.28
$ws.Range("N58").Value = -6399.316

$ws.Range("H62").Value = 60895.39
$ws.Range("J62").Value = 6819.636
$ws.Range("L62").Value = 6819.636
$ws.Range("N62").Value = -8067.636

$ws.Range("H65").Value = 60895.39
$ws.Range("J65").Value = 6819.636
$ws.Range("L65").Value = 34098.18
$ws.Range("N65").Value = -40338.18

$ws.Range("H99").Value = 8811.280000000001
$ws.Range("I99").Value = 5024.1875
$ws.Range("K99").Value = 5024.1875
$ws.Range("M99").Value = -3526.1875

$ws.Range("H101").Value = 13861.2
$ws.Range("J101").Value = 13861.2
$ws.Range("L101").Value = 13861.2
$ws.Range("N101").Value = -20351.2

$ws.Range("H107").Value = 864.7692
$ws.Range("I107").Value = 893.5
$ws.Range("J107").Value = 769
$ws.Range("K107").Value = 893.5
$ws.Range("L107").Value = 769
$ws.Range("M107").Value = 1026.5
$ws.Range("N107").Value = -4609

$ws.Range("H117").Value = 90000
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H126").Value = 8811.280000000001
$ws.Range("I126").Value = 5024.1875
$ws.Range("K126").Value = 15072.5625
$ws.Range("M126").Value = -12602.5625

$ws.Range("H134").Value = 5806.3228
$ws.Range("I134").Value = 4464.8887
$ws.Range("J134").Value = 7663.6924
$ws.Range("K134").Value = 13394.6661
$ws.Range("L134").Value = 22991.0772
$ws.Range("M134").Value = -10859.6661
$ws.Range("N134").Value = -28061.0772

$ws.Range("H135").Value = 65000
$ws.Range("J135").Value = 65000
$ws.Range("L135").Value = 65000
$ws.Range("N135").Value = -75140

$ws.Range("H136").Value = 3838.75
$ws.Range("I136").Value = 2201.28
$ws.Range("J136").Value = 5993.316
$ws.Range("K136").Value = 6603.84
$ws.Range("L136").Value = 17979.948
$ws.Range("M136").Value = -4053.84
$ws.Range("N136").Value = -23079.948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1090.1052
$ws.Range("I12").Value = 32.666668
$ws.Range("J12").Value = 1288.375
$ws.Range("K12").Value = 98.000004
$ws.Range("L12").Value = 3865.125
$ws.Range("M12").Value = 74.999996
$ws.Range("N12").Value = -4211.125

$ws.Range("H45").Value = 2054.1
$ws.Range("J45").Value = 1393.4445
$ws.Range("L45").Value = 4180.333500000001
$ws.Range("N45").Value = -5244.333500000001

$ws.Range("H92").Value = 2323
$ws.Range("I92").Value = 1689
$ws.Range("K92").Value = 5067
$ws.Range("M92").Value = -3819

$ws.Range("H102").Value = 11026
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H107").Value = 1195.2084
$ws.Range("I107").Value = 1278.2858
$ws.Range("K107").Value = 3834.8574
$ws.Range("M107").Value = -1914.8574

$ws.Range("H124").Value = 2069.375
$ws.Range("I124").Value = 3685
$ws.Range("J124").Value = 1100
$ws.Range("K124").Value = 11055
$ws.Range("L124").Value = 3300
$ws.Range("M124").Value = -6145
$ws.Range("N124").Value = -13120

$ws.Range("H130").Value = 2566.6667
$ws.Range("J130").Value = 3000
$ws.Range("L130").Value = 9000
$ws.Range("N130").Value = -19040

$ws.Range("H141").Value = 2971.7778
$ws.Range("I141").Value = 2971.7778
$ws.Range("K141").Value = 8915.3334
$ws.Range("M141").Value = -3735.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 279042.5
$ws.Range("I62").Value = 58000
$ws.Range("K62").Value = 58000
$ws.Range("M62").Value = -57314

$ws.Range("H65").Value = 279042.5
$ws.Range("I65").Value = 58000
$ws.Range("K65").Value = 174000
$ws.Range("M65").Value = -170568

$ws.Range("H102").Value = 1590.6562
$ws.Range("I102").Value = 1613.4
$ws.Range("K102").Value = 1613.4
$ws.Range("M102").Value = 8.599999999999909

$ws.Range("H113").Value = 88889.87
$ws.Range("I113").Value = 119834.35
$ws.Range("K113").Value = 119834.35
$ws.Range("M113").Value = -117664.35

$ws.Range("H122").Value = 1626
$ws.Range("I122").Value = 1618.6154
$ws.Range("K122").Value = 4855.8462
$ws.Range("M122").Value = -2405.8462

$ws.Range("H132").Value = 17073.867
$ws.Range("I132").Value = 19474.666
$ws.Range("K132").Value = 58423.99800000001
$ws.Range("M132").Value = -55893.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 846.25
$ws.Range("I22").Value = 995
$ws.Range("J22").Value = 796.6667
$ws.Range("K22").Value = 995
$ws.Range("L22").Value = 796.6667
$ws.Range("M22").Value = -700
$ws.Range("N22").Value = -1386.6667

$ws.Range("H27").Value = 846.25
$ws.Range("I27").Value = 995
$ws.Range("J27").Value = 796.6667
$ws.Range("K27").Value = 995
$ws.Range("L27").Value = 796.6667
$ws.Range("M27").Value = -888
$ws.Range("N27").Value = -1010.6667

$ws.Range("H32").Value = 3873.75
$ws.Range("I32").Value = 1833.3334
$ws.Range("K32").Value = 1833.3334
$ws.Range("M32").Value = -1516.3334

$ws.Range("H61").Value = 2071.389
$ws.Range("I61").Value = 2071.389
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2071.389
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1869.389
$ws.Range("N61").ClearContents()

$ws.Range("H93").Value = 6062.7144
$ws.Range("I93").Value = 5831.5625
$ws.Range("J93").Value = 6802.4
$ws.Range("K93").Value = 5831.5625
$ws.Range("L93").Value = 6802.4
$ws.Range("M93").Value = -4583.5625
$ws.Range("N93").Value = -9298.4

$ws.Range("H99").Value = 28991.777
$ws.Range("I99").Value = 28991.777
$ws.Range("K99").Value = 28991.777
$ws.Range("M99").Value = -25996.777

$ws.Range("H113").Value = 2071.389
$ws.Range("I113").Value = 2071.389
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2071.389
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 98.61099999999988
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 3827.889
$ws.Range("I122").Value = 3182.611
$ws.Range("J122").Value = 5118.4443
$ws.Range("K122").Value = 9547.832999999999
$ws.Range("L122").Value = 15355.3329
$ws.Range("M122").Value = -7097.832999999999
$ws.Range("N122").Value = -20255.3329

$ws.Range("H132").Value = 5109.9736
$ws.Range("I132").Value = 4858.9395
$ws.Range("K132").Value = 14576.8185
$ws.Range("M132").Value = -12046.8185

$ws.Range("H136").Value = 7076.2915
$ws.Range("I136").Value = 6481.294
$ws.Range("J136").Value = 8521.286
$ws.Range("K136").Value = 19443.882
$ws.Range("L136").Value = 25563.858
$ws.Range("M136").Value = -16893.882
$ws.Range("N136").Value = -30663.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 9592.111000000001
$ws.Range("I26").Value = 9904.143
$ws.Range("K26").Value = 9904.143
$ws.Range("M26").Value = -9611.143

$ws.Range("H54").Value = 22638.643
$ws.Range("J54").Value = 24911.75
$ws.Range("L54").Value = 24911.75
$ws.Range("N54").Value = -25951.75

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

$ws.Range("H100").Value = 1519.909
$ws.Range("I100").Value = 1648
$ws.Range("K100").Value = 3296
$ws.Range("M100").Value = -2755

$ws.Range("H119").Value = 265316.66
$ws.Range("J119").Value = 265316.66
$ws.Range("L119").Value = 265316.66
$ws.Range("N119").Value = -274992.66

$ws.Range("H122").Value = 1578.3334
$ws.Range("I122").Value = 1287.2858
$ws.Range("J122").Value = 1985.8
$ws.Range("K122").Value = 3861.8574
$ws.Range("L122").Value = 5957.4
$ws.Range("M122").Value = -1411.8574
$ws.Range("N122").Value = -10857.4

$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800

$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840

$ws.Range("H126").Value = 5226.8237
$ws.Range("I126").Value = 3407.459
$ws.Range("K126").Value = 10222.377
$ws.Range("M126").Value = -7752.377

$ws.Range("H132").Value = 138521.06
$ws.Range("I132").Value = 177144.27
$ws.Range("K132").Value = 531432.8099999999
$ws.Range("M132").Value = -528902.8099999999

$ws.Range("H136").Value = 4879793.5
$ws.Range("I136").Value = 6452830.5
$ws.Range("J136").Value = 3377.5
$ws.Range("K136").Value = 19358491.5
$ws.Range("L136").Value = 10132.5
$ws.Range("M136").Value = -19355941.5
$ws.Range("N136").Value = -15232.5

Write-Host "done"
